$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.014927392184047891
$ws.Range("B1").Value = -0.014927392923491314

$ws.Range("A2").Value = 0.014151785833180178
$ws.Range("B2").Value = -0.014151786588835676

$ws.Range("A3").Value = -0.048182502341397329
$ws.Range("B3").Value = 0.048182501606491318

$ws.Range("A4").Value = 0.075191258073497289
$ws.Range("B4").Value = -0.075191258761216823
